$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two extra exported columns.
$ws.Range("J1").Value = "nextkin"
$ws.Range("K1").Value = "kinphone"

# Match the header formatting (bold / fill / border) used by the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1:K1").PasteSpecial(-4122) | Out-Null

# Match the plain bordered formatting used by the rest of the data rows.
$ws.Range("H2:H14").Copy() | Out-Null
$ws.Range("J2:K14").PasteSpecial(-4122) | Out-Null

# Widen the new columns the same as column H.
$w = $ws.Columns.Item(8).ColumnWidth
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 10)).EntireColumn.ColumnWidth = $w

# Extend the "duplicateValues" conditional formatting that covered H2:H14 to also
# cover the two new columns (H2:J14).
$cf = $ws.Range("H2:H14").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("H2:J14"))

# Restore the active selection reported by the author's session.
$ws.Range("J8").Select() | Out-Null

$excel.CutCopyMode = $false
